$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '39.389.84'
$ws.Range("E2").Value = '  +1.65%  '
$ws.Range("D3").Value = '2.157.43'
$ws.Range("E3").Value = '  +3.12%  '
$ws.Range("D5").Value = '''227.75'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.50%  '
$ws.Range("D6").Value = '''0.623'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.96%  '
$ws.Range("D7").Value = '''64.13'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +5.01%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").Value = '''0.397'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.88%  '
$ws.Range("D10").Value = '''0.0859'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.73%  '
$ws.Range("E11").Value = '  +0.01%  '
$ws.Range("D12").Value = '''15.99'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.40%  '
$ws.Range("D13").Value = '2.477.22'
$ws.Range("E13").Value = '  +3.10%  '
$ws.Range("D14").Value = '''22.31'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.24%  '
$ws.Range("D15").Value = '''0.812'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.80%  '
$ws.Range("E16").Value = '  +0.87%  '
$ws.Range("D17").Value = '2.157.95'
$ws.Range("E17").Value = '  +3.27%  '
$ws.Range("D18").Value = '39.340.55'
$ws.Range("E18").Value = '  +1.68%  '
$ws.Range("D19").Value = '''71.79'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.06%  '
$ws.Range("E20").Value = '  +0.44%  '
$ws.Range("D21").Value = '0.0₃0854'
$ws.Range("E21").Value = '  +1.53%  '
$ws.Range("D22").Value = '''231.29'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.52%  '
$ws.Range("D24").Value = '''2.51'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.94%  '
$ws.Range("D25").Value = '''2.35'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.57%  '
$ws.Range("D26").Value = '''172.07'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.40%  '
$ws.Range("E27").Value = '  -0.09%  '
$ws.Range("E28").Value = '  +1.04%  '
$ws.Range("D29").Value = '''19.91'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.01%  '
$ws.Range("E30").Value = '  -1.46%  '
$ws.Range("E31").Value = '  +8.83%  '
$ws.Range("E32").Value = '  +1.09%  '
$ws.Range("D33").Value = '''4.61'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.10%  '
$ws.Range("D34").Value = '''4.75'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.24%  '
$ws.Range("D35").Value = '''7.07'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +8.88%  '
$ws.Range("D36").Value = '''0.0619'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.59%  '
$ws.Range("D37").Value = '''2.40'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.34%  '
$ws.Range("E38").Value = '  -0.22%  '
$ws.Range("E39").Value = '  +0.05%  '
$ws.Range("D40").Value = '''103.86'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.81%  '
$ws.Range("E41").Value = '  +0.98%  '
$ws.Range("D42").Value = '''17.81'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.50%  '
$ws.Range("D43").Value = '1.539.83'
$ws.Range("E43").Value = '  +0.26%  '
$ws.Range("E44").Value = '  +4.23%  '
$ws.Range("B45").Value = 'HuobiToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D45").Value = '''2.82'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.58%  '
$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").Value = '''0.0925'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.62%  '
$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").Value = '''7.80'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.83%  '
$ws.Range("E48").Value = '  +5.54%  '
$ws.Range("D49").Value = '''4.23'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.98%  '
$ws.Range("D50").Value = '2.360.63'
$ws.Range("E50").Value = '  +3.14%  '
$ws.Range("E51").Value = '  +0.13%  '
